$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column (H) - header, matching style of existing header cells (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Populate data values for the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
